# Trading update: 2026-02-18 10:58:14
# Renames the "TestStrategy" placeholder into the real "MarketMaking"
# strategy, records two new live trades for it, refreshes the "Strategy
# Status" board with the full current roster, and rolls the Summary
# capital/active-strategy counters forward.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the 4th sheet (and the strategy it tracks) from the
#    temporary "TestStrategy" name to "MarketMaking".
# ---------------------------------------------------------------------
$wsStrategy = $wb.Worksheets.Item("TestStrategy")
$wsStrategy.Name = "MarketMaking"

# ---------------------------------------------------------------------
# 2) Summary sheet: capital was funded and the strategy roster is live.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(2, 2).Value = 1500   # Initial Capital
$wsSummary.Cells.Item(3, 2).Value = 1500   # Current Capital
$wsSummary.Cells.Item(11, 2).Value = 15    # Active Strategies

# ---------------------------------------------------------------------
# 3) Strategy Status sheet: populate the full roster of active
#    strategies (currently only the header row exists).
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")

$strategies = @(
    "EMAArbitrage",
    "HighProbConvergence",
    "HighProbabilityBond",
    "LongshotBias",
    "MarketMaking",
    "MicrostructureScalper",
    "arbitrage",
    "breakout_momentum",
    "leadlag",
    "momentum",
    "orderbook_imbalance",
    "sentiment",
    "sharp_money",
    "volatility_scorer",
    "vwap"
)

$row = 2
foreach ($name in $strategies) {
    $wsStatus.Cells.Item($row, 1).Value = $name
    $wsStatus.Cells.Item($row, 2).Value = "ACTIVE"
    $wsStatus.Cells.Item($row, 3).Value = 100
    $wsStatus.Cells.Item($row, 4).Value = 0
    $wsStatus.Cells.Item($row, 5).Value = 0
    $wsStatus.Cells.Item($row, 6).Value = 0
    $wsStatus.Cells.Item($row, 7).Value = 0
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4) All Trades sheet: trade #7's exit/confidence/entry-reason columns
#    get cleared back out (still open, no fill yet), and two new
#    MarketMaking trades (#8 and #9) get appended.
# ---------------------------------------------------------------------
$wsTrades = $wb.Worksheets.Item("All Trades")

# Trade #7 (row 8): exit price now tracked as a live 0 quote, and the
# capital/slippage/confidence/entry-reason columns go back to blank.
$wsTrades.Cells.Item(8, 7).Value = 0
$wsTrades.Cells.Item(8, 11).ClearContents()
$wsTrades.Cells.Item(8, 12).ClearContents()
$wsTrades.Cells.Item(8, 13).ClearContents()
$wsTrades.Cells.Item(8, 14).ClearContents()
$wsTrades.Cells.Item(8, 15).ClearContents()
$wsTrades.Cells.Item(8, 16).ClearContents()
$wsTrades.Cells.Item(8, 17).ClearContents()

# Trade #8 (row 9): new MarketMaking trade, still open.
# (the leading "'" keeps the YYYY-MM-DD text from being auto-converted
# into a date serial by the COM layer - it is stripped from the stored
# value, same as typing it into Excel directly.)
$wsTrades.Cells.Item(9, 1).Value = 8
$wsTrades.Cells.Item(9, 2).Value = "'2026-02-18"
$wsTrades.Cells.Item(9, 3).Value = "10:56:25"
$wsTrades.Cells.Item(9, 4).Value = "MarketMaking"
$wsTrades.Cells.Item(9, 5).Value = "DOWN"
$wsTrades.Cells.Item(9, 6).Value = 0.47
$wsTrades.Cells.Item(9, 7).Value = 0
$wsTrades.Cells.Item(9, 8).Value = "OPEN"
$wsTrades.Cells.Item(9, 9).Value = 0
$wsTrades.Cells.Item(9, 10).Value = 0
$wsTrades.Cells.Item(9, 11).ClearContents()
$wsTrades.Cells.Item(9, 12).ClearContents()
$wsTrades.Cells.Item(9, 13).ClearContents()
$wsTrades.Cells.Item(9, 14).ClearContents()
$wsTrades.Cells.Item(9, 15).ClearContents()
$wsTrades.Cells.Item(9, 16).ClearContents()
$wsTrades.Cells.Item(9, 17).ClearContents()

# Trade #9 (row 10): new MarketMaking trade, still open, with a filled
# confidence/entry-reason pair (this is also the latest trade, mirrored
# onto the MarketMaking sheet below).
$wsTrades.Cells.Item(10, 1).Value = 9
$wsTrades.Cells.Item(10, 2).Value = "'2026-02-18"
$wsTrades.Cells.Item(10, 3).Value = "10:57:43"
$wsTrades.Cells.Item(10, 4).Value = "MarketMaking"
$wsTrades.Cells.Item(10, 5).Value = "DOWN"
$wsTrades.Cells.Item(10, 6).Value = 0.68
$wsTrades.Cells.Item(10, 7).ClearContents()
$wsTrades.Cells.Item(10, 8).Value = "OPEN"
$wsTrades.Cells.Item(10, 9).Value = 0
$wsTrades.Cells.Item(10, 10).Value = 0
$wsTrades.Cells.Item(10, 11).Value = 100
$wsTrades.Cells.Item(10, 12).Value = 0
$wsTrades.Cells.Item(10, 13).Value = 0
$wsTrades.Cells.Item(10, 14).Value = 0.6
$wsTrades.Cells.Item(10, 15).Value = "Normal spread capture: 202 bps"
$wsTrades.Cells.Item(10, 16).ClearContents()
$wsTrades.Cells.Item(10, 17).Value = 0

# ---------------------------------------------------------------------
# 5) MarketMaking sheet (the old "TestStrategy" snapshot row): replace
#    the stale test-entry row with the latest real trade (#9), matching
#    row 10 of "All Trades" above.
# ---------------------------------------------------------------------
$wsStrategy.Cells.Item(2, 1).Value = 9
$wsStrategy.Cells.Item(2, 3).Value = "10:57:43"
$wsStrategy.Cells.Item(2, 4).Value = "MarketMaking"
$wsStrategy.Cells.Item(2, 5).Value = "DOWN"
$wsStrategy.Cells.Item(2, 6).Value = 0.68
$wsStrategy.Cells.Item(2, 14).Value = 0.6
$wsStrategy.Cells.Item(2, 15).Value = "Normal spread capture: 202 bps"
